$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 85.71429000000001
$ws.Range("I39").Value = 80
$ws.Range("J39").Value = 93.333336
$ws.Range("K39").Value = 240
$ws.Range("L39").Value = 280.000008
$ws.Range("M39").Value = 56
$ws.Range("N39").Value = -872.000008

# Row 64
$ws.Range("H64").Value = 287145.53
$ws.Range("I64").Value = 321542.28
$ws.Range("J64").Value = 3372.5
$ws.Range("K64").Value = 321542.28
$ws.Range("L64").Value = 3372.5
$ws.Range("M64").Value = -321294.28
$ws.Range("N64").Value = -3868.5

# Row 67
$ws.Range("H67").Value = 287145.53
$ws.Range("I67").Value = 321542.28
$ws.Range("J67").Value = 3372.5
$ws.Range("K67").Value = 321542.28
$ws.Range("L67").Value = 3372.5
$ws.Range("M67").Value = -320684.28
$ws.Range("N67").Value = -5088.5

# Row 113
$ws.Range("H113").Value = 2218.6843
$ws.Range("J113").Value = 2738.889
$ws.Range("L113").Value = 2738.889
$ws.Range("N113").Value = -9246.888999999999

# Row 116
$ws.Range("H116").Value = 14112.223
$ws.Range("I116").Value = 100005
$ws.Range("J116").Value = 3375.625
$ws.Range("K116").Value = 100005
$ws.Range("L116").Value = 3375.625
$ws.Range("M116").Value = -96563
$ws.Range("N116").Value = -10259.625

# Row 132
$ws.Range("H132").Value = 26318646
$ws.Range("I132").Value = 33335320
$ws.Range("J132").Value = 6125.125
$ws.Range("K132").Value = 100005960
$ws.Range("L132").Value = 18375.375
$ws.Range("M132").Value = -100003430
$ws.Range("N132").Value = -23435.375

# Row 137
$ws.Range("H137").Value = 870.1389
$ws.Range("I137").Value = 778
$ws.Range("J137").Value = 943.85
$ws.Range("K137").Value = 2334
$ws.Range("L137").Value = 2831.55
$ws.Range("M137").Value = 216
$ws.Range("N137").Value = -7931.55

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 360019.62
$ws.Range("I32").Value = 2422.908
$ws.Range("K32").Value = 2422.908
$ws.Range("M32").Value = -2135.908

# Row 61
$ws.Range("H61").Value = 1074.0444
$ws.Range("I61").Value = 736.9091
$ws.Range("J61").Value = 2001.1666
$ws.Range("K61").Value = 736.9091
$ws.Range("L61").Value = 2001.1666
$ws.Range("M61").Value = -524.9091
$ws.Range("N61").Value = -2425.1666

# Row 102
$ws.Range("H102").Value = 3224.88
$ws.Range("I102").Value = 2505.2856
$ws.Range("J102").Value = 7002.75
$ws.Range("K102").Value = 2505.2856
$ws.Range("L102").Value = 7002.75
$ws.Range("M102").Value = -883.2856000000002
$ws.Range("N102").Value = -10246.75

# Row 132
$ws.Range("H132").Value = 1391.1562
$ws.Range("I132").Value = 971.5599999999999
$ws.Range("K132").Value = 2914.68
$ws.Range("M132").Value = -384.6799999999998

# Row 136
$ws.Range("H136").Value = 1074.0444
$ws.Range("I136").Value = 736.9091
$ws.Range("J136").Value = 2001.1666
$ws.Range("K136").Value = 2210.7273
$ws.Range("L136").Value = 6003.4998
$ws.Range("M136").Value = 339.2727
$ws.Range("N136").Value = -11103.4998

$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 28000
$ws.Range("J62").Value = 28000
$ws.Range("L62").Value = 28000
$ws.Range("N62").Value = -29372

# Row 65
$ws.Range("H65").Value = 28000
$ws.Range("J65").Value = 28000
$ws.Range("L65").Value = 84000
$ws.Range("N65").Value = -90864

# Row 99
$ws.Range("H99").Value = 1539.1364
$ws.Range("I99").Value = 842.3077
$ws.Range("J99").Value = 2545.6667
$ws.Range("K99").Value = 842.3077
$ws.Range("L99").Value = 2545.6667
$ws.Range("M99").Value = 655.6923
$ws.Range("N99").Value = -5541.6667

# Row 134
$ws.Range("H134").Value = 3598.4814
$ws.Range("I134").Value = 651.7568
$ws.Range("J134").Value = 10011.941
$ws.Range("K134").Value = 1955.2704
$ws.Range("L134").Value = 30035.823
$ws.Range("M134").Value = 579.7296000000001
$ws.Range("N134").Value = -35105.823

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 1514.1904
$ws.Range("I99").Value = 1510.421
$ws.Range("J99").Value = 1550
$ws.Range("K99").Value = 1510.421
$ws.Range("L99").Value = 1550
$ws.Range("M99").Value = -12.42100000000005
$ws.Range("N99").Value = -4546

# Row 126
$ws.Range("H126").Value = 1514.1904
$ws.Range("I126").Value = 1510.421
$ws.Range("J126").Value = 1550
$ws.Range("K126").Value = 4531.263
$ws.Range("L126").Value = 4650
$ws.Range("M126").Value = -2061.263
$ws.Range("N126").Value = -9590

# Row 132
$ws.Range("H132").Value = 29619.629
$ws.Range("I132").Value = 977.88
$ws.Range("J132").Value = 101224
$ws.Range("K132").Value = 2933.64
$ws.Range("L132").Value = 303672
$ws.Range("M132").Value = -403.6399999999999
$ws.Range("N132").Value = -308732

# Row 134
$ws.Range("H134").Value = 24987.62
$ws.Range("I134").Value = 27507.37
$ws.Range("J134").Value = 1050
$ws.Range("K134").Value = 82522.11
$ws.Range("L134").Value = 3150
$ws.Range("M134").Value = -79987.11
$ws.Range("N134").Value = -8220

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 6330067.5
$ws.Range("I131").Value = 814.53845
$ws.Range("J131").Value = 7576738.5
$ws.Range("K131").Value = 2443.61535
$ws.Range("L131").Value = 22730215.5
$ws.Range("M131").Value = 2596.38465
$ws.Range("N131").Value = -22740295.5

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 35347.066
$ws.Range("I132").Value = 41942.96
$ws.Range("J132").Value = 2367.6
$ws.Range("K132").Value = 125828.88
$ws.Range("L132").Value = 7102.799999999999
$ws.Range("M132").Value = -123298.88
$ws.Range("N132").Value = -12162.8

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1635.1428
$ws.Range("I82").Value = 1593.8948
$ws.Range("J82").Value = 1722.2222
$ws.Range("K82").Value = 1593.8948
$ws.Range("L82").Value = 1722.2222
$ws.Range("M82").Value = -1232.8948
$ws.Range("N82").Value = -2444.2222

# Row 85
$ws.Range("H85").Value = 1635.1428
$ws.Range("I85").Value = 1593.8948
$ws.Range("J85").Value = 1722.2222
$ws.Range("K85").Value = 1593.8948
$ws.Range("L85").Value = 1722.2222
$ws.Range("M85").Value = -345.8948
$ws.Range("N85").Value = -4218.2222

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 87358.336
$ws.Range("I96").Value = 2133.3333
$ws.Range("J96").Value = 115766.664
$ws.Range("K96").Value = 2133.3333
$ws.Range("L96").Value = 115766.664
$ws.Range("M96").Value = -760.3332999999998
$ws.Range("N96").Value = -118512.664

# Row 132
$ws.Range("H132").Value = 16027915
$ws.Range("I132").Value = 21930836
$ws.Range("J132").Value = 5702.381
$ws.Range("K132").Value = 65792508
$ws.Range("L132").Value = 17107.143
$ws.Range("M132").Value = -65789978
$ws.Range("N132").Value = -22167.143

# Row 136
$ws.Range("H136").Value = 2651.2708
$ws.Range("I136").Value = 3426.8484
$ws.Range("J136").Value = 945
$ws.Range("K136").Value = 10280.5452
$ws.Range("L136").Value = 2835
$ws.Range("M136").Value = -7730.5452
$ws.Range("N136").Value = -7935
